# The workbook tracks daily Coliflor (cauliflower) price records for the
# "Macroferia Regional de Talca" market. A new daily record was inserted
# into the weekly series at row 195, pushing the existing rows 195-268
# down to 196-269 (dimension grows from A1:R268 to A1:R269).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 195; this shifts rows 195-268 down
# to 196-269 automatically (and updates the sheet dimension accordingly).
$ws.Rows(195).Insert()

# Populate the newly inserted row 195 with the new record's data.
$ws.Range("A195").Value = 5
$ws.Range("B195").Value = "Macroferia Regional de Talca"
$ws.Range("C195").Value = "Maule"
$ws.Range("D195").Value = 44755
$ws.Range("E195").Value = 7
$ws.Range("F195").Value = 100112008
$ws.Range("G195").Value = "Coliflor"
$ws.Range("H195").Value = "Sin especificar"
$ws.Range("I195").Value = "Primera"
$ws.Range("J195").Value = 3000
$ws.Range("K195").Value = 1000
$ws.Range("L195").Value = 1000
$ws.Range("M195").Value = 1000
$ws.Range("N195").Value = "`$/unidad"
$ws.Range("O195").Value = "Región del Maule"
$ws.Range("P195").Value = 1000
$ws.Range("Q195").Value = 1
$ws.Range("R195").Value = "Hortaliza"
